# Applies the "output folder feature" edit:
#  - adds a source-filename column to the Input sheet
#  - adds two "Monte Carlo error success rate" columns to the Calc sheet
#  - updates a handful of recalculated numeric results on Calc/Results/Constants
#  - widens/narrows a few columns to fit the new content

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: Excel's ColumnWidth setter (in this runtime) only lands on
# multiples of 1/6 of a character, offset by 5/6. Pick the closest reachable
# input for a desired stored width so the resulting file is as close as
# possible to the target column width.
# ---------------------------------------------------------------------------
function Set-ClosestColumnWidth($col, $targetWidth) {
    $offset = 0.8333333333333334
    $k = [math]::Round(($targetWidth - $offset) * 6)
    $input = $k / 6.0
    $col.ColumnWidth = $input
}

# ===========================================================================
# Sheet "Input"
# ===========================================================================
$wsInput = $wb.Worksheets.Item("Input")

# Column D now holds the long "Filename: ..." text -> widen it
Set-ClosestColumnWidth $wsInput.Columns.Item(4) 66.7109375

# New cell D3: source filename for this analysis
$wsInput.Range("D3").Value = "Filename: C:\Neptune\User\Neptune\Data\UTh\2015\0815\011_7184.dat"

# ===========================================================================
# Sheet "Calc"
# ===========================================================================
$wsCalc = $wb.Worksheets.Item("Calc")

# Column width tweaks to accommodate new/resized columns
Set-ClosestColumnWidth $wsCalc.Columns.Item(43) 20.7109375
Set-ClosestColumnWidth $wsCalc.Columns.Item(49) 9.7109375
Set-ClosestColumnWidth $wsCalc.Columns.Item(51) 20.7109375
Set-ClosestColumnWidth $wsCalc.Columns.Item(58) 18.7109375
Set-ClosestColumnWidth $wsCalc.Columns.Item(59) 32.7109375
Set-ClosestColumnWidth $wsCalc.Columns.Item(60) 30.7109375

# New columns BG/BH: "Unkorr./Korr. Montefehler Erfolgsrate", formatted like BF
$wsCalc.Range("BF1:BF2").Copy($wsCalc.Range("BG1:BG2"))
$wsCalc.Range("BF1:BF2").Copy($wsCalc.Range("BH1:BH2"))
$wsCalc.Range("BG1").Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Range("BH1").Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BG2").Value = "(%)"
$wsCalc.Range("BH2").Value = "(%)"

# Recalculated numeric results (row 3)
$wsCalc.Range("AP3").Value = 0.5246
$wsCalc.Range("AQ3").Value = 0.1857624386649325
$wsCalc.Range("AW3").Value = 0.5295
$wsCalc.Range("AX3").Value = 0.5342210465052091
$wsCalc.Range("AY3").Value = 0.187532605704462
$wsCalc.Range("BC3").Value = 0.568838764942721
$wsCalc.Range("BE3").Value = 267.1105232526045
$wsCalc.Range("BF3").Value = 0.189204655095914

# New Monte-Carlo success-rate values
$wsCalc.Range("BG3").Value = 100
$wsCalc.Range("BH3").Value = 100

# ===========================================================================
# Sheet "Results"
# ===========================================================================
$wsResults = $wb.Worksheets.Item("Results")

Set-ClosestColumnWidth $wsResults.Columns.Item(16) 8.7109375

$wsResults.Range("N3").Value = 0.5246
$wsResults.Range("P3").Value = 0.5295
$wsResults.Range("R3").Value = 0.568838764942721

# ===========================================================================
# Sheet "Constants"
# ===========================================================================
$wsConstants = $wb.Worksheets.Item("Constants")

$wsConstants.Range("B3").Value = 0.00005
